$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "no_effect" column before the existing weather-effect columns (C)
$ws.Columns("C").Insert()

# Insert a new "midi" column inside the mc1 cue/monitor block (now at Q after first insert)
$ws.Columns("Q").Insert()

# Populate the new "midi" column (Q) first so its shared-string slot is allocated
# before the "no_effect" one, matching the append order of the source workbook.
$ws.Range("Q1").Value = "midi"
$ws.Range("Q2").Value = "903a7f"

# Populate the new "no_effect" column (C)
$ws.Range("C1").Value = "no_effect"
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"

# Update the current selection to match the edited workbook state
$null = $ws.Range("C6").Select()
